{"js": "// COA Assignment 1 - text corrections\n// 1) \"Non restoring\" -> \"Non-restoring\" (hyphenate)\n// 2) \"input /output\" -> \"input/output\" (drop stray space before the slash)\n// 3) Re-type the \"...explain Distributed Arbitration scheme.\" sentence right after\n//    \"With the help of an example\" (instead of after the _GoBack bookmark) and add a\n//    trailing space, which leaves the _GoBack bookmark sitting at the end of the\n//    paragraph - mirroring how Word re-anchors _GoBack at the point of the last edit.\n\nconst body = context.document.body;\n\n// --- 1) Hyphenate \"Non restoring\" -> \"Non-restoring\" ---------------------------------\nconst hyphenTargets = body.search(\"Non restoring\", { matchCase: true });\nhyphenTargets.load(\"items\");\nawait context.sync();\nfor (const r of hyphenTargets.items) {\n  r.insertText(\"Non-restoring\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) Remove the stray space in \"input /output\" -------------------------------------\nconst slashTargets = body.search(\"input /output\", { matchCase: true });\nslashTargets.load(\"items\");\nawait context.sync();\nfor (const r of slashTargets.items) {\n  r.insertText(\"input/output\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 3) Move the \"explain Distributed Arbitration scheme.\" clause ahead of the bookmark\nconst clauseTargets = body.search(\", explain Distributed Arbitration scheme.\", { matchCase: true });\nclauseTargets.load(\"items\");\nawait context.sync();\nfor (const r of clauseTargets.items) {\n  r.insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst exampleTargets = body.search(\"With the help of an example\", { matchCase: true });\nexampleTargets.load(\"items\");\nawait context.sync();\nconst exampleEnd = exampleTargets.items[0].getRange(Word.RangeLocation.end);\nexampleEnd.insertText(\", explain Distributed Arbitration scheme. \", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# COA Assignment 1 - text corrections\n# 1) \"Non restoring\" -> \"Non-restoring\" (hyphenate)\n# 2) \"input /output\" -> \"input/output\" (drop stray space before the slash)\n# 3) Re-type the \"...explain Distributed Arbitration scheme.\" sentence right after\n#    \"With the help of an example\" (instead of after the _GoBack bookmark) and add a\n#    trailing space, which leaves the _GoBack bookmark sitting at the end of the\n#    paragraph - mirroring how Word re-anchors _GoBack at the point of the last edit.\n\n$d = $word.ActiveDocument\n\n# --- 1) Hyphenate \"Non restoring\" -> \"Non-restoring\" ---------------------------------\n$range1 = $d.Content\n$range1.Find.Execute(\"Non restoring\", $false, $false, $false, $false, $false, $true, 1, $false, \"Non-restoring\", 2)\n\n# --- 2) Remove the stray space in \"input /output\" -------------------------------------\n$range2 = $d.Content\n$range2.Find.Execute(\"input /output\", $false, $false, $false, $false, $false, $true, 1, $false, \"input/output\", 2)\n\n# --- 3) Move the \"explain Distributed Arbitration scheme.\" clause ahead of the bookmark\n$range3 = $d.Content\n$range3.Find.Execute(\", explain Distributed Arbitration scheme.\")\n$range3.Text = \"\"\n\n$range4 = $d.Content\n$range4.Find.Execute(\"With the help of an example\")\n$range4.Collapse(0)\n$range4.InsertAfter(\", explain Distributed Arbitration scheme. \")\n"}
